# Add 2022-Q3 data:
#  - insert a new "2022-Q3" sheet (copied from "2022-Q2" so header/index styling
#    matches the other quarter sheets) right after "总计" and before "2022-Q2"
#  - populate it with the 2022-Q3 fund-holding data
#  - update the "总计" (summary) sheet with the new quarter's totals row
#  - fix the "2021-Q3" sheet's D1 header text ("基金金额" -> "基金规模")

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet by copying "2022-Q2" (keeps the same
#    header / index-column cell styles as the other quarterly sheets) and
#    placing it immediately before "2022-Q2".
# ---------------------------------------------------------------------------
$sourceSheet = $wb.Worksheets.Item("2022-Q2")
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$sourceSheet.Copy($beforeSheet)

$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2. Populate "2022-Q3" with the new data (6 funds, one more row than the
#    previous quarter, so the existing row 6 style is reused for row 7).
# ---------------------------------------------------------------------------
$q3Data = @(
    @("159869", "华夏中证动漫游戏ETF",       "6.35", "99.31", "4.33", "0.2750", 8),
    @("516010", "国泰中证动漫游戏ETF",       "3.78", "97.86", "4.18", "0.1580", 8),
    @("161030", "富国中证体育产业指数A",     "1.59", "94.00", "4.17", "0.0663", 9),
    @("516770", "华泰柏瑞中证动漫游戏ETF",   "0.99", "96.39", "4.30", "0.0426", 8),
    @("013278", "富国中证体育产业指数C",     "0.42", "94.00", "4.17", "0.0175", 9),
    @("517500", "国泰中证沪港深动漫游戏ETF", "0.53", "92.78", "3.03", "0.0161", 10)
)

for ($i = 0; $i -lt $q3Data.Length; $i++) {
    $row = $i + 2
    $rec = $q3Data[$i]
    $q3.Cells.Item($row, 1).Value = $i
    $q3.Cells.Item($row, 2).Value = $rec[0]
    $q3.Cells.Item($row, 3).Value = $rec[1]
    $q3.Cells.Item($row, 4).Value = $rec[2]
    $q3.Cells.Item($row, 5).Value = $rec[3]
    $q3.Cells.Item($row, 6).Value = $rec[4]
    $q3.Cells.Item($row, 7).Value = $rec[5]
    $q3.Cells.Item($row, 8).Value = $rec[6]
}

# ---------------------------------------------------------------------------
# 3. Update the "总计" summary sheet: shift the existing quarters down one
#    row and add the new 2022-Q3 total at the top of the data.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summaryData = @(
    @("2022-Q3", 6, 0.58),
    @("2022-Q2", 5, 0.62),
    @("2022-Q1", 4, 0.64),
    @("2021-Q4", 4, 0.64),
    @("2021-Q3", 4, 0.52),
    @("2021-Q2", 5, 0.59)
)

for ($i = 0; $i -lt $summaryData.Length; $i++) {
    $row = $i + 2
    $rec = $summaryData[$i]
    $summary.Cells.Item($row, 1).Value = $i
    $summary.Cells.Item($row, 2).Value = $rec[0]
    $summary.Cells.Item($row, 3).Value = $rec[1]
    $summary.Cells.Item($row, 4).Value = $rec[2]
}

# ---------------------------------------------------------------------------
# 4. Fix the "2021-Q3" sheet header: "基金金额" -> "基金规模".
# ---------------------------------------------------------------------------
$q3_2021 = $wb.Worksheets.Item("2021-Q3")
$q3_2021.Range("D1").Value = "基金规模"
